$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "2007年" (row 2). This shifts the remaining
# data rows (2010年, 2012年, 2015年) up by one, matching the diff.
$ws.Rows.Item(2).Delete()
